# Auto-generated Excel COM-interop script to apply scheduled-runner price/profit updates
# across the 8 Leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 135.625
$ws.Range("I33").Value = 147.85715
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 147.85715
$ws.Range("L33").Value = 50
$ws.Range("M33").Value = 81.14285000000001
$ws.Range("N33").Value = -508

$ws.Range("H107").Value = 1252.1818
$ws.Range("I107").Value = 586.1111
$ws.Range("K107").Value = 586.1111
$ws.Range("M107").Value = 1333.8889

$ws.Range("H111").Value = 3474.2942
$ws.Range("I111").Value = 2183.6667
$ws.Range("J111").Value = 4926.25
$ws.Range("K111").Value = 6551.000100000001
$ws.Range("L111").Value = 14778.75
$ws.Range("M111").Value = -3484.000100000001
$ws.Range("N111").Value = -20912.75

$ws.Range("H138").Value = 3364.8823
$ws.Range("I138").Value = 3214.6667
$ws.Range("J138").Value = 3533.875
$ws.Range("K138").Value = 9644.000100000001
$ws.Range("L138").Value = 10601.625
$ws.Range("M138").Value = -4504.000100000001
$ws.Range("N138").Value = -20881.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2113.6667
$ws.Range("I2").Value = 2036.4
$ws.Range("K2").Value = 2036.4
$ws.Range("M2").Value = -1923.4

$ws.Range("H32").Value = 17697
$ws.Range("I32").Value = 15399.8
$ws.Range("K32").Value = 15399.8
$ws.Range("M32").Value = -15112.8

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H110").Value = 945.2222
$ws.Range("I110").Value = 945.2222
$ws.Range("K110").Value = 945.2222
$ws.Range("M110").Value = 1099.7778

$ws.Range("H116").Value = 2113.6667
$ws.Range("I116").Value = 2036.4
$ws.Range("K116").Value = 2036.4
$ws.Range("M116").Value = 257.5999999999999

$ws.Range("H122").Value = 3499.6667
$ws.Range("I122").Value = 3499.6667
$ws.Range("K122").Value = 10499.0001
$ws.Range("M122").Value = -8049.000100000001

$ws.Range("H132").Value = 1554.8
$ws.Range("I132").Value = 1471.9231
$ws.Range("J132").Value = 2093.5
$ws.Range("K132").Value = 4415.7693
$ws.Range("L132").Value = 6280.5
$ws.Range("M132").Value = -1885.7693
$ws.Range("N132").Value = -11340.5

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2113.6667
$ws.Range("I3").Value = 2036.4
$ws.Range("K3").Value = 2036.4
$ws.Range("M3").Value = -1922.4

$ws.Range("H94").Value = 2163.8572
$ws.Range("I94").Value = 2191.1667
$ws.Range("K94").Value = 2191.1667
$ws.Range("M94").Value = -1740.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 8047.3335
$ws.Range("I39").Value = 1499.9375
$ws.Range("J39").Value = 28999
$ws.Range("K39").Value = 1499.9375
$ws.Range("L39").Value = 28999
$ws.Range("M39").Value = -1108.9375
$ws.Range("N39").Value = -29781

$ws.Range("H49").Value = 8047.3335
$ws.Range("I49").Value = 1499.9375
$ws.Range("J49").Value = 28999
$ws.Range("K49").Value = 1499.9375
$ws.Range("L49").Value = 28999
$ws.Range("M49").Value = -1317.9375
$ws.Range("N49").Value = -29363

$ws.Range("H62").Value = 2833.3333
$ws.Range("J62").Value = 2833.3333
$ws.Range("L62").Value = 2833.3333
$ws.Range("N62").Value = -4081.3333

$ws.Range("H65").Value = 2833.3333
$ws.Range("J65").Value = 2833.3333
$ws.Range("L65").Value = 14166.6665
$ws.Range("N65").Value = -20406.6665

$ws.Range("H132").Value = 2926.8333
$ws.Range("I132").Value = 3265.25
$ws.Range("K132").Value = 9795.75
$ws.Range("M132").Value = -7265.75

$ws.Range("H134").Value = 1283.3334
$ws.Range("I134").Value = 1283.3334
$ws.Range("K134").Value = 3850.0002
$ws.Range("M134").Value = -1315.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 10
$ws.Range("I22").Value = 10
$ws.Range("K22").Value = 30
$ws.Range("M22").Value = 139

$ws.Range("H27").Value = 10
$ws.Range("I27").Value = 10
$ws.Range("K27").Value = 30
$ws.Range("M27").Value = 72

$ws.Range("H35").Value = 7787.75
$ws.Range("I35").Value = 10075
$ws.Range("K35").Value = 30225
$ws.Range("M35").Value = -29937

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H129").Value = 1820.3
$ws.Range("I129").Value = 2228
$ws.Range("J129").Value = 1208.75
$ws.Range("K129").Value = 6684
$ws.Range("L129").Value = 3626.25
$ws.Range("M129").Value = -1684
$ws.Range("N129").Value = -13626.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3756.5
$ws.Range("I102").Value = 3683.353
$ws.Range("K102").Value = 3683.353
$ws.Range("M102").Value = -2061.353

$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960

$ws.Range("H130").Value = 55000
$ws.Range("J130").Value = 55000
$ws.Range("L130").Value = 55000
$ws.Range("N130").Value = -65040

$ws.Range("H132").Value = 2563.9167
$ws.Range("I132").Value = 2263
$ws.Range("J132").Value = 3466.6667
$ws.Range("K132").Value = 6789
$ws.Range("L132").Value = 10400.0001
$ws.Range("M132").Value = -4259
$ws.Range("N132").Value = -15460.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6666.6665
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 7500
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = -4888
$ws.Range("N7").Value = -7724

$ws.Range("H126").Value = 6666.6665
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -27440

$ws.Range("H132").Value = 22495.092
$ws.Range("I132").Value = 22555.75
$ws.Range("K132").Value = 67667.25
$ws.Range("M132").Value = -65137.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 16000
$ws.Range("I15").Value = 12000
$ws.Range("K15").Value = 12000
$ws.Range("M15").Value = -11712

$ws.Range("H113").Value = 559.2
$ws.Range("I113").Value = 510.22223
$ws.Range("K113").Value = 1530.66669
$ws.Range("M113").Value = 639.33331

$ws.Range("H126").Value = 1747.3
$ws.Range("I126").Value = 882
$ws.Range("K126").Value = 2646
$ws.Range("M126").Value = -176

$ws.Range("H132").Value = 8285
$ws.Range("I132").Value = 5056.5
$ws.Range("K132").Value = 15169.5
$ws.Range("M132").Value = -12639.5

Write-Output "Applied scheduled-runner updates to 37 rows across 8 sheets."